$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos table values. Force text format first so numeric-looking
# strings (prices, percentages) are stored verbatim instead of being
# auto-converted to floating point numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.282.32'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.529.72'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.84%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.60'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.62'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.522.09'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.76%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.87'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.385'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.137.26'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000182'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.07'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.545.65'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.03%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.260.07'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.13'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.90'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.22'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +5.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '391.47'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.572'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.681.29'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.75'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.99%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000114'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +9.19%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +11.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.15'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.551.27'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.21%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.79'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.144'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +16.00%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '169.75'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.56'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +7.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.91'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.99'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +6.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0800'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.823'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.44'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +19.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.58'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.41'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.68'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.19'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +9.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.78'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.392.52'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +10.30%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '304.04'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +10.43%  '
